$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-09-22"

# Update the "September (through 09-21)" label cell to "September (through 09-22)"
$ws.Range("A10").Value = "September (through 09-22)"

# Update September row (row 10) values for 2017-2022 (columns D-I)
$ws.Range("D10").Value = 51
$ws.Range("E10").Value = 42
$ws.Range("F10").Value = 55
$ws.Range("G10").Value = 84
$ws.Range("H10").Value = 132
$ws.Range("I10").Value = 107

# Update Total row (row 11) values for 2017-2022 (columns D-I)
$ws.Range("D11").Value = 602
$ws.Range("E11").Value = 532
$ws.Range("F11").Value = 404
$ws.Range("G11").Value = 868
$ws.Range("H11").Value = 1202
$ws.Range("I11").Value = 1242
